$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A166").Value = 45983
$ws.Range("B166").Value = "四方坪站充电量(kw)"
$ws.Range("C166").Value = 670.0139999999999
$ws.Range("D166").Value = 923.12
$ws.Range("E166").Value = 483.62
$ws.Range("F166").Value = 398.62099999999998
$ws.Range("G166").Value = 300.78399999999999
$ws.Range("H166").Value = 838.82900000000006
$ws.Range("I166").Value = 408.012
$ws.Range("J166").Value = 141.018
$ws.Range("K166").Value = 96.125
$ws.Range("L166").Value = 121.24600000000001
$ws.Range("M166").Value = 199.822
$ws.Range("N166").Value = 235.10000000000002
$ws.Range("O166").Value = 556.5089999999999
$ws.Range("P166").Value = 1625.6069999999997
$ws.Range("Q166").Value = 646.00100000000009
$ws.Range("R166").Value = 349.05099999999993
$ws.Range("S166").Value = 499.56
$ws.Range("T166").Value = 150.54599999999999
$ws.Range("U166").Value = 141.13000000000002
$ws.Range("V166").Value = 255.21000000000004
$ws.Range("W166").Value = 215.99199999999999
$ws.Range("X166").Value = 60.317999999999998
$ws.Range("Y166").Value = 93.5
$ws.Range("Z166").Value = 94

$ws.Range("A167").Value = 45983
$ws.Range("B167").Value = "高岭站充电量(kw)"
$ws.Range("C167").Value = 624.41999999999985
$ws.Range("D167").Value = 459.88699999999994
$ws.Range("E167").Value = 288.83699999999999
$ws.Range("F167").Value = 156.98099999999999
$ws.Range("G167").Value = 98.022999999999996
$ws.Range("H167").Value = 57.707999999999998
$ws.Range("I167").Value = 346.673
$ws.Range("J167").Value = 165.02100000000002
$ws.Range("K167").Value = 367.81400000000002
$ws.Range("L167").Value = 295.80100000000004
$ws.Range("M167").Value = 188.79500000000002
$ws.Range("N167").Value = 165.91199999999998
$ws.Range("O167").Value = 633.11400000000003
$ws.Range("P167").Value = 191.10000000000002
$ws.Range("Q167").Value = 226.989
$ws.Range("R167").Value = 143.38900000000001
$ws.Range("S167").Value = 202.471
$ws.Range("T167").Value = 131.398
$ws.Range("U167").Value = 115.98400000000001
$ws.Range("V167").Value = 51.144000000000005
$ws.Range("W167").Value = 200.77900000000002
$ws.Range("X167").Value = 9.7810000000000006
$ws.Range("Y167").Value = 33.226999999999997
$ws.Range("Z167").Value = 59.389000000000003

$ws.Range("A168").Value = 45984
$ws.Range("B168").Value = "四方坪站充电量(kw)"
$ws.Range("C168").Value = 662.43500000000006
$ws.Range("D168").Value = 1006.067
$ws.Range("E168").Value = 531.57600000000002
$ws.Range("F168").Value = 454.58300000000003
$ws.Range("G168").Value = 278.31899999999996
$ws.Range("H168").Value = 557.27199999999993
$ws.Range("I168").Value = 368.34300000000007
$ws.Range("J168").Value = 56.472999999999999
$ws.Range("K168").Value = 148.601
$ws.Range("L168").Value = 61.36
$ws.Range("M168").Value = 187.15199999999999
$ws.Range("N168").Value = 131.00799999999998
$ws.Range("O168").Value = 488.03399999999999
$ws.Range("P168").Value = 1463.778
$ws.Range("Q168").Value = 825.5200000000001
$ws.Range("R168").Value = 643.52200000000016
$ws.Range("S168").Value = 289.87099999999998
$ws.Range("T168").Value = 198.47100000000003
$ws.Range("U168").Value = 49.7
$ws.Range("V168").Value = 142.62
$ws.Range("W168").Value = 74.911000000000001
$ws.Range("X168").Value = 39.54
$ws.Range("Y168").Value = 37.200000000000003
$ws.Range("Z168").Value = 0

$ws.Range("A169").Value = 45984
$ws.Range("B169").Value = "高岭站充电量(kw)"
$ws.Range("C169").Value = 318.24700000000007
$ws.Range("D169").Value = 365.11699999999996
$ws.Range("E169").Value = 195.81199999999998
$ws.Range("F169").Value = 67.710999999999999
$ws.Range("G169").Value = 39.772999999999996
$ws.Range("H169").Value = 90.694999999999993
$ws.Range("I169").Value = 247.95500000000001
$ws.Range("J169").Value = 258.05400000000003
$ws.Range("K169").Value = 151.75899999999999
$ws.Range("L169").Value = 104.02799999999999
$ws.Range("M169").Value = 191.88900000000001
$ws.Range("N169").Value = 244.84200000000001
$ws.Range("O169").Value = 296.13199999999995
$ws.Range("P169").Value = 441.30899999999997
$ws.Range("Q169").Value = 64.082999999999998
$ws.Range("R169").Value = 256.76099999999997
$ws.Range("S169").Value = 158.41799999999998
$ws.Range("T169").Value = 37.664000000000001
$ws.Range("U169").Value = 72.281000000000006
$ws.Range("V169").Value = 13.064
$ws.Range("W169").Value = 26.445
$ws.Range("X169").Value = 40.381999999999998
$ws.Range("Y169").Value = 58.048000000000002
$ws.Range("Z169").Value = 30.913999999999998

$ws.Range("A170").Value = 45985
$ws.Range("B170").Value = "四方坪站充电量(kw)"
$ws.Range("C170").Value = 850.55300000000022
$ws.Range("D170").Value = 752.85500000000002
$ws.Range("E170").Value = 250.8
$ws.Range("F170").Value = 352.55799999999994
$ws.Range("G170").Value = 328.435
$ws.Range("H170").Value = 381.35400000000004
$ws.Range("I170").Value = 504.73099999999999
$ws.Range("J170").Value = 166.95100000000002
$ws.Range("K170").Value = 154.411
$ws.Range("L170").Value = 218.387
$ws.Range("M170").Value = 121.7
$ws.Range("N170").Value = 209.42400000000001
$ws.Range("O170").Value = 823.70200000000023
$ws.Range("P170").Value = 1409.3730000000005
$ws.Range("Q170").Value = 321.92399999999998
$ws.Range("R170").Value = 475.02100000000007
$ws.Range("S170").Value = 409.90600000000001
$ws.Range("T170").Value = 264.06600000000003
$ws.Range("U170").Value = 44.597999999999999
$ws.Range("V170").Value = 107.04
$ws.Range("W170").Value = 111.16000000000001
$ws.Range("X170").Value = 140.346
$ws.Range("Y170").Value = 55.26
$ws.Range("Z170").Value = 80.323000000000008

$ws.Range("A171").Value = 45985
$ws.Range("B171").Value = "高岭站充电量(kw)"
$ws.Range("C171").Value = 694.00099999999998
$ws.Range("D171").Value = 354.22600000000006
$ws.Range("E171").Value = 126.453
$ws.Range("F171").Value = 48.798999999999999
$ws.Range("G171").Value = 25.234000000000002
$ws.Range("H171").Value = 109.613
$ws.Range("I171").Value = 78.415000000000006
$ws.Range("J171").Value = 99.015999999999991
$ws.Range("K171").Value = 224.22400000000002
$ws.Range("L171").Value = 149.459
$ws.Range("M171").Value = 388.99599999999998
$ws.Range("N171").Value = 238.41600000000003
$ws.Range("O171").Value = 368.13499999999999
$ws.Range("P171").Value = 422.96100000000001
$ws.Range("Q171").Value = 274.47699999999998
$ws.Range("R171").Value = 60.598999999999997
$ws.Range("S171").Value = 16.227
$ws.Range("T171").Value = 35.423999999999999
$ws.Range("U171").Value = 5.6719999999999997
$ws.Range("V171").Value = 10.516
$ws.Range("W171").Value = 43.164000000000001
$ws.Range("X171").Value = 0
$ws.Range("Y171").Value = 0
$ws.Range("Z171").Value = 4.07


$ws.Range("E174").Select() | Out-Null

